$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names (A2:A5) - replace old part names with new appliance names
$ws.Range("A2").Value = "lavadora"
$ws.Range("A3").Value = "correa"
$ws.Range("A4").Value = "bujias"
$ws.Range("A5").Value = "lamparita"

# Update prices (B2:B5)
$ws.Range("B2").Value = 400
$ws.Range("B3").Value = 700
$ws.Range("B4").Value = 300
$ws.Range("B5").Value = 100

# Update the active selection to G1
$ws.Range("G1").Select() | Out-Null
